$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; this shifts columns B:F left to A:E
$ws.Range("A1:A7").EntireColumn.Delete()
